$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1891891891891892
$ws.Range("C2").Value = 0.5598455598455598
$ws.Range("J2").Value = 0.007722007722007722
$ws.Range("P2").Value = 0.1003861003861004
$ws.Range("C3").Value = 0.006756756756756757
$ws.Range("J3").Value = 0.02027027027027027
$ws.Range("P3").Value = 0.7702702702702703
$ws.Range("S3").Value = 0.2027027027027027
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.2972972972972973
$ws.Range("B6").Value = 0.1059322033898305
$ws.Range("D6").Value = 0.00423728813559322
$ws.Range("F6").Value = 0.0635593220338983
$ws.Range("J6").Value = 0.2372881355932203
$ws.Range("O6").Value = 0.01694915254237288
$ws.Range("Q6").Value = 0.1822033898305085
$ws.Range("R6").Value = 0.0635593220338983
$ws.Range("S6").Value = 0.326271186440678
$ws.Range("B7").Value = 0.04205607476635514
$ws.Range("D7").Value = 0.03271028037383177
$ws.Range("F7").Value = 0.06542056074766354
$ws.Range("J7").Value = 0.1401869158878505
$ws.Range("O7").Value = 0.009345794392523364
$ws.Range("Q7").Value = 0.2289719626168224
$ws.Range("R7").Value = 0.0514018691588785
$ws.Range("S7").Value = 0.4299065420560748
$ws.Range("B8").Value = 0.08658008658008658
$ws.Range("D8").Value = 0.01731601731601732
$ws.Range("F8").Value = 0.06277056277056277
$ws.Range("J8").Value = 0.08441558441558442
$ws.Range("O8").Value = 0.02164502164502164
$ws.Range("Q8").Value = 0.1991341991341991
$ws.Range("R8").Value = 0.08658008658008658
$ws.Range("S8").Value = 0.4415584415584415
$ws.Range("B9").Value = 0.0776255707762557
$ws.Range("D9").Value = 0.0136986301369863
$ws.Range("F9").Value = 0.0502283105022831
$ws.Range("J9").Value = 0.091324200913242
$ws.Range("O9").Value = 0.0182648401826484
$ws.Range("Q9").Value = 0.1917808219178082
$ws.Range("R9").Value = 0.1050228310502283
$ws.Range("S9").Value = 0.4520547945205479
$ws.Range("B10").Value = 0.09915254237288136
$ws.Range("D10").Value = 0.0211864406779661
$ws.Range("E10").Value = 0.000847457627118644
$ws.Range("F10").Value = 0.08559322033898305
$ws.Range("J10").Value = 0.08983050847457627
$ws.Range("O10").Value = 0.02288135593220339
$ws.Range("Q10").Value = 0.2161016949152542
$ws.Range("R10").Value = 0.08050847457627118
$ws.Range("S10").Value = 0.3838983050847458
$ws.Range("G11").Value = 0.0932475884244373
$ws.Range("J11").Value = 0.1061093247588424
$ws.Range("K11").Value = 0.1639871382636656
$ws.Range("L11").Value = 0.6237942122186495
$ws.Range("S11").Value = 0.01286173633440514
$ws.Range("G12").Value = 0.7864077669902912
$ws.Range("J12").Value = 0.1650485436893204
$ws.Range("L12").Value = 0.02912621359223301
$ws.Range("S12").Value = 0.01941747572815534
$ws.Range("G13").Value = 0.6226415094339622
$ws.Range("J13").Value = 0.3018867924528302
$ws.Range("S13").Value = 0.07547169811320754
$ws.Range("F15").Value = 0.02083333333333333
$ws.Range("H15").Value = 0.1833333333333333
$ws.Range("I15").Value = 0.07083333333333333
$ws.Range("J15").Value = 0.3208333333333334
$ws.Range("K15").Value = 0.08333333333333333
$ws.Range("O15").Value = 0.05833333333333333
$ws.Range("S15").Value = 0.2625
$ws.Range("H16").Value = 0.1635220125786163
$ws.Range("I16").Value = 0.07547169811320754
$ws.Range("J16").Value = 0.4465408805031447
$ws.Range("K16").Value = 0.119496855345912
$ws.Range("M16").Value = 0.03144654088050314
$ws.Range("O16").Value = 0.05660377358490566
$ws.Range("S16").Value = 0.1069182389937107
$ws.Range("F17").Value = 0.02536997885835095
$ws.Range("H17").Value = 0.1670190274841438
$ws.Range("I17").Value = 0.105708245243129
$ws.Range("J17").Value = 0.4059196617336152
$ws.Range("K17").Value = 0.105708245243129
$ws.Range("M17").Value = 0.02959830866807611
$ws.Range("O17").Value = 0.07610993657505286
$ws.Range("S17").Value = 0.08456659619450317
$ws.Range("F18").Value = 0.01092896174863388
$ws.Range("H18").Value = 0.2131147540983606
$ws.Range("I18").Value = 0.1038251366120219
$ws.Range("J18").Value = 0.3879781420765027
$ws.Range("K18").Value = 0.09836065573770492
$ws.Range("M18").Value = 0.02185792349726776
$ws.Range("O18").Value = 0.06010928961748634
$ws.Range("S18").Value = 0.1038251366120219
$ws.Range("F19").Value = 0.01932367149758454
$ws.Range("H19").Value = 0.2230273752012882
$ws.Range("I19").Value = 0.0966183574879227
$ws.Range("J19").Value = 0.3599033816425121
$ws.Range("K19").Value = 0.1183574879227053
$ws.Range("M19").Value = 0.0249597423510467
$ws.Range("O19").Value = 0.07326892109500806
$ws.Range("S19").Value = 0.08454106280193237
